# Fix typo on slide 3 ("03-cvicenie"): "zozmiera" -> "zozbiera"
# (opravene zadanie na zvicenie cislo 3)

$p = $ppt.ActivePresentation

$needle     = "zozmiera"
$replacement = "zozbiera"
$fixed = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $tf = $shape.TextFrame
        if (-not $tf.HasText) {
            continue
        }

        $tr = $tf.TextRange
        $text = $tr.Text
        $idx = $text.IndexOf($needle)

        if ($idx -ge 0) {
            # Target only the exact run of characters that misspell the word,
            # leaving the rest of the run/paragraph formatting untouched.
            $run = $tr.Characters($idx + 1, $needle.Length)
            $run.Text = $replacement
            $fixed = $true
        }
    }
}

if (-not $fixed) {
    throw "Could not find text '$needle' to fix in the presentation."
}
